$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws 'D2' '20.196.13'
Set-TextValue $ws 'E2' '  +1.30%  '
Set-TextValue $ws 'D3' '1.442.01'
Set-TextValue $ws 'E3' '  +2.48%  '
Set-TextValue $ws 'E4' '  +0.63%  '
Set-TextValue $ws 'D5' '0.9174'
Set-TextValue $ws 'E5' '  -8.37%  '
Set-TextValue $ws 'D6' '277.23'
Set-TextValue $ws 'E6' '  +2.11%  '
Set-TextValue $ws 'D7' '0.3667'
Set-TextValue $ws 'E7' '  -0.71%  '
Set-TextValue $ws 'D8' '0.3135'
Set-TextValue $ws 'E8' '  +2.44%  '
Set-TextValue $ws 'D9' '38.97'
Set-TextValue $ws 'E9' '  -0.06%  '
Set-TextValue $ws 'E10' '  +4.50%  '
Set-TextValue $ws 'D11' '0.06538'
Set-TextValue $ws 'E11' '  +0.52%  '
Set-TextValue $ws 'D12' '1.002'
Set-TextValue $ws 'E12' '  -0.01%  '
Set-TextValue $ws 'D13' '5.407'
Set-TextValue $ws 'E13' '  +1.73%  '
Set-TextValue $ws 'D14' '17.62'
Set-TextValue $ws 'E14' '  +5.04%  '
Set-TextValue $ws 'D15' '6.081'
Set-TextValue $ws 'E15' '  -0.34%  '
Set-TextValue $ws 'D16' '1.443.84'
Set-TextValue $ws 'E16' '  +2.59%  '
Set-TextValue $ws 'D17' '0.00001016'
Set-TextValue $ws 'E17' '  +1.23%  '
Set-TextValue $ws 'D18' '0.9369'
Set-TextValue $ws 'E18' '  -6.43%  '
Set-TextValue $ws 'D19' '0.05630'
Set-TextValue $ws 'E19' '  -1.35%  '
Set-TextValue $ws 'D20' '67.59'
Set-TextValue $ws 'E20' '  -6.78%  '
Set-TextValue $ws 'D21' '5.433'
Set-TextValue $ws 'E21' '  -2.05%  '
Set-TextValue $ws 'D22' '14.51'
Set-TextValue $ws 'E22' '  +1.83%  '
Set-TextValue $ws 'D23' '10.85'
Set-TextValue $ws 'E23' '  +0.82%  '
Set-TextValue $ws 'D24' '2.268'
Set-TextValue $ws 'E24' '  -0.14%  '
Set-TextValue $ws 'D25' '20.206.49'
Set-TextValue $ws 'E25' '  +1.25%  '
Set-TextValue $ws 'D26' '2.192'
Set-TextValue $ws 'D27' '136.41'
Set-TextValue $ws 'E27' '  -0.68%  '
Set-TextValue $ws 'D28' '16.97'
Set-TextValue $ws 'E28' '  +2.31%  '
Set-TextValue $ws 'D29' '1.601.61'
Set-TextValue $ws 'E29' '  +2.34%  '
Set-TextValue $ws 'D30' '111.09'
Set-TextValue $ws 'E30' '  +2.81%  '
Set-TextValue $ws 'D31' '3.773'
Set-TextValue $ws 'E31' '  -1.75%  '
Set-TextValue $ws 'D32' '0.8159'
Set-TextValue $ws 'E32' '  +0.41%  '
Set-TextValue $ws 'D33' '4.844'
Set-TextValue $ws 'E33' '  -7.22%  '
Set-TextValue $ws 'D34' '0.07693'
Set-TextValue $ws 'E34' '  +0.27%  '
Set-TextValue $ws 'B35' 'Hedera'
Set-TextValue $ws 'C35' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D35' '0.06037'
Set-TextValue $ws 'E35' '  +4.01%  '
Set-TextValue $ws 'B36' 'WEMIXTOKEN'
Set-TextValue $ws 'C36' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D36' '1.487'
Set-TextValue $ws 'E36' '  +15.24%  '
Set-TextValue $ws 'E37' '  -0.86%  '
Set-TextValue $ws 'D38' '1.131'
Set-TextValue $ws 'E38' '  +6.29%  '
Set-TextValue $ws 'D39' '10.26'
Set-TextValue $ws 'E40' '  -1.63%  '
Set-TextValue $ws 'D41' '0.9355'
Set-TextValue $ws 'E41' '  -6.54%  '
Set-TextValue $ws 'D42' '0.1830'
Set-TextValue $ws 'E42' '  -5.67%  '
Set-TextValue $ws 'D43' '7.137'
Set-TextValue $ws 'E43' '  -14.56%  '
Set-TextValue $ws 'D44' '0.5252'
Set-TextValue $ws 'E44' '  -0.06%  '
Set-TextValue $ws 'D45' '3.520'
Set-TextValue $ws 'E45' '  +0.38%  '
Set-TextValue $ws 'D46' '12.07'
Set-TextValue $ws 'E46' '  -0.07%  '
Set-TextValue $ws 'D47' '119.56'
Set-TextValue $ws 'E47' '  +8.96%  '
Set-TextValue $ws 'E48' '  +1.67%  '
Set-TextValue $ws 'D49' '1.773'
Set-TextValue $ws 'E49' '  -0.15%  '
Set-TextValue $ws 'D50' '0.06337'
Set-TextValue $ws 'E50' '  +3.01%  '
Set-TextValue $ws 'D51' '0.9950'
Set-TextValue $ws 'E51' '  -0.64%  '
